$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Body: insert a new "Order custom PCBs" checklist item right before the
#    existing "Order hardware components" item.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd() -eq "Order hardware components") {
        $para.Range.InsertParagraphBefore()
        $d.Paragraphs($i).Range.Text = "Order custom PCBs"
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Body: rename the "Print 'User Quick Guide'" checklist item to
#    "Print 'User Guide'".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Print " + [char]0x201C + "User Quick Guide" + [char]0x201D,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Print " + [char]0x201C + "User Guide" + [char]0x201D, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Body: rename the stand-alone "'User Quick Guide'" reference (in the
#    "Items to Give to User" section) to "'User Guide'".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    [char]0x201C + "User Quick Guide" + [char]0x201D,
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x201C + "User Guide" + [char]0x201D, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Header: merge the "OpenAT" + "-Switch-Latch" runs (and drop the
#    now-orphaned spell-check proofErr markers) into a single
#    "OpenAT-Switch-Latch" run.
# ---------------------------------------------------------------------------
$hf = $d.Sections(1).Headers(1)
for ($i = 1; $i -le $hf.Range.Paragraphs.Count; $i++) {
    $p = $hf.Range.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "OpenAT-Switch-Latch") {
        $p.Range.InsertParagraphBefore()
        $hf.Range.Paragraphs($i).Range.Text = "OpenAT-Switch-Latch"
        $hf.Range.Paragraphs($i + 1).Range.Delete()
        break
    }
}

Write-Output "done"
